$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "stock" info column (H) values
$ws.Range("H17").Value = "en stock"
$ws.Range("H19").Value = "pas nécessaire"

# Restore selection to B10, matching the authored session state
$ws.Range("B10").Select()
